{"js": "// Add a new row at the end of the (only) table in the document body,\n// right after the \"Walter Murch\" / \"Filmmaking\" row, with the new\n// publication entry described in the commit diff.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRowValues = [\n  [\"10\", \"Baltimoore\", \"\u041f\u0435\u0447\u0430\u0442\u043d\u044b\u0439\", \"fdsafadfdafd\", \"N/A\", \"N.Tasbolatuly, A.Bekzhan\"]\n];\n\ntable.addRows(\"End\", 1, newRowValues);\nawait context.sync();\n", "ps1": "# Add a new row at the end of the (only) table in the document, right\n# after the \"Walter Murch\" / \"Filmmaking\" row, with the new publication\n# entry described in the commit diff.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"10\"\n$newRow.Cells.Item(2).Range.Text = \"Baltimoore\"\n$newRow.Cells.Item(3).Range.Text = \"\u041f\u0435\u0447\u0430\u0442\u043d\u044b\u0439\"\n$newRow.Cells.Item(4).Range.Text = \"fdsafadfdafd\"\n$newRow.Cells.Item(5).Range.Text = \"N/A\"\n$newRow.Cells.Item(6).Range.Text = \"N.Tasbolatuly, A.Bekzhan\"\n"}
